$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new column before C (shifts old C..K to D..L), making room for
# a new "Thai ID" data column right after "Seq No".
$ws.Columns("C:C").Insert()

# The inserted column should be the same width as column B.
$ws.Columns("C:C").ColumnWidth = 11.785714285714286

# The "To" label (originally in C4) must stay in column C -- move it back
# from D4 (where the column insert pushed it) and clear the now-empty D4.
$ws.Range("D4").Cut($ws.Range("C4"))
$ws.Range("D4").Clear()

# Update header row text: "SeqNo" -> "Seq No", and the new column gets
# the new "Thai ID" header.
$ws.Range("B7").Value = "Seq No"
$ws.Range("C7").Value = "Thai ID"

# Grow the hidden _xlnm._FilterDatabase range by one column to include
# the new column.
$wb.Names.Item(1).RefersTo = "=Sheet1!`$A`$7:`$H`$7"

# Restore the active cell/selection.
$ws.Activate()
$ws.Range("E10").Select()
